# Auto-generated Excel COM-interop script to apply the diff changes
# Updates LeveProfits calculation columns (H-N) across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 632.25
$ws.Range("I2").Value = 195.5
$ws.Range("J2").Value = 850.625
$ws.Range("K2").Value = 195.5
$ws.Range("L2").Value = 850.625
$ws.Range("M2").Value = -82.5
$ws.Range("N2").Value = -1076.625

$ws.Range("H9").Value = 174.75
$ws.Range("I9").Value = 170
$ws.Range("J9").Value = 176.33333
$ws.Range("K9").Value = 170
$ws.Range("L9").Value = 176.33333
$ws.Range("M9").Value = -1
$ws.Range("N9").Value = -514.3333299999999

$ws.Range("H12").Value = 188
$ws.Range("I12").Value = 172
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 172
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -2
$ws.Range("N12").Value = -640

$ws.Range("H19").Value = 1584.75
$ws.Range("I19").Value = 853.625
$ws.Range("J19").Value = 1877.2
$ws.Range("K19").Value = 853.625
$ws.Range("L19").Value = 1877.2
$ws.Range("M19").Value = -678.625
$ws.Range("N19").Value = -2227.2

$ws.Range("H33").Value = 162.35
$ws.Range("I33").Value = 103.35714
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 103.35714
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = 125.64286
$ws.Range("N33").Value = -758

$ws.Range("H38").Value = 985.6667
$ws.Range("I38").Value = 182.8
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 548.4000000000001
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -176.4000000000001
$ws.Range("N38").Value = -15744

$ws.Range("H58").Value = 7875
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 7875
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 23625
$ws.Range("N58").Value = -23925

$ws.Range("H61").Value = 15
$ws.Range("I61").Value = 15
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 45
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 127

$ws.Range("H70").Value = 2118
$ws.Range("I70").Value = 3210
$ws.Range("J70").Value = 1162.5
$ws.Range("K70").Value = 9630
$ws.Range("L70").Value = 3487.5
$ws.Range("M70").Value = -9360
$ws.Range("N70").Value = -4027.5

$ws.Range("H73").Value = 2118
$ws.Range("I73").Value = 3210
$ws.Range("J73").Value = 1162.5
$ws.Range("K73").Value = 9630
$ws.Range("L73").Value = 3487.5
$ws.Range("M73").Value = -8694
$ws.Range("N73").Value = -5359.5

$ws.Range("H137").Value = 1586.1428
$ws.Range("I137").Value = 1328.6364
$ws.Range("J137").Value = 2530.3333
$ws.Range("K137").Value = 3985.9092
$ws.Range("L137").Value = 7590.999899999999
$ws.Range("M137").Value = -1435.9092
$ws.Range("N137").Value = -12690.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -1884

$ws.Range("H5").Value = 97.5
$ws.Range("I5").Value = 88.333336
$ws.Range("J5").Value = 125
$ws.Range("K5").Value = 88.333336
$ws.Range("L5").Value = 125
$ws.Range("M5").Value = 23.666664
$ws.Range("N5").Value = -349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 97.5
$ws.Range("I4").Value = 88.333336
$ws.Range("J4").Value = 125
$ws.Range("K4").Value = 88.333336
$ws.Range("L4").Value = 125
$ws.Range("M4").Value = 26.666664
$ws.Range("N4").Value = -355

$ws.Range("H22").Value = 1205.9
$ws.Range("I22").Value = 196.5
$ws.Range("J22").Value = 2720
$ws.Range("K22").Value = 196.5
$ws.Range("L22").Value = 2720
$ws.Range("M22").Value = -23.5
$ws.Range("N22").Value = -3066

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 156.88235
$ws.Range("I7").Value = 129.75
$ws.Range("J7").Value = 222
$ws.Range("K7").Value = 129.75
$ws.Range("L7").Value = 222
$ws.Range("M7").Value = -16.75
$ws.Range("N7").Value = -448

$ws.Range("H16").Value = 50003116
$ws.Range("I16").Value = 55558264
$ws.Range("J16").Value = 6799
$ws.Range("K16").Value = 55558264
$ws.Range("L16").Value = 6799
$ws.Range("M16").Value = -55557977
$ws.Range("N16").Value = -7373

$ws.Range("H22").Value = 91807.414
$ws.Range("I22").Value = 91062.55
$ws.Range("J22").Value = 100001
$ws.Range("K22").Value = 91062.55
$ws.Range("L22").Value = 100001
$ws.Range("M22").Value = -90712.55
$ws.Range("N22").Value = -100701

$ws.Range("H25").Value = 2073.5557
$ws.Range("I25").Value = 1616.6666
$ws.Range("J25").Value = 2439.0667
$ws.Range("K25").Value = 1616.6666
$ws.Range("L25").Value = 2439.0667
$ws.Range("M25").Value = -1442.6666
$ws.Range("N25").Value = -2787.0667

$ws.Range("H31").Value = 4071.2856
$ws.Range("I31").Value = 2747.8462
$ws.Range("J31").Value = 6221.875
$ws.Range("K31").Value = 2747.8462
$ws.Range("L31").Value = 6221.875
$ws.Range("M31").Value = -2452.8462
$ws.Range("N31").Value = -6811.875

$ws.Range("H34").Value = 4071.2856
$ws.Range("I34").Value = 2747.8462
$ws.Range("J34").Value = 6221.875
$ws.Range("K34").Value = 2747.8462
$ws.Range("L34").Value = 6221.875
$ws.Range("M34").Value = -2545.8462
$ws.Range("N34").Value = -6625.875

$ws.Range("H113").Value = 50003116
$ws.Range("I113").Value = 55558264
$ws.Range("J113").Value = 6799
$ws.Range("K113").Value = 55558264
$ws.Range("L113").Value = 6799
$ws.Range("M113").Value = -55556094
$ws.Range("N113").Value = -11139

$ws.Range("H132").Value = 2767.3333
$ws.Range("I132").Value = 2459.1667
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7377.500100000001
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4847.500100000001
$ws.Range("N132").Value = -17060

$ws.Range("H134").Value = 1876.7
$ws.Range("I134").Value = 1502.1765
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 4506.529500000001
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -1971.529500000001
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 243.5
$ws.Range("I2").Value = 40.625
$ws.Range("J2").Value = 514
$ws.Range("K2").Value = 243.75
$ws.Range("L2").Value = 3084
$ws.Range("M2").Value = -130.75
$ws.Range("N2").Value = -3310

$ws.Range("H7").Value = 20000138
$ws.Range("I7").Value = 20000138
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 60000414
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -60000302

$ws.Range("H17").Value = 365
$ws.Range("I17").Value = 45
$ws.Range("J17").Value = 525
$ws.Range("K17").Value = 135
$ws.Range("L17").Value = 1575
$ws.Range("M17").Value = 34
$ws.Range("N17").Value = -1913

$ws.Range("H19").Value = 9500
$ws.Range("I19").Value = 10000
$ws.Range("J19").Value = 9000
$ws.Range("K19").Value = 30000
$ws.Range("L19").Value = 27000
$ws.Range("M19").Value = -29826
$ws.Range("N19").Value = -27348

$ws.Range("H23").Value = 250147.25
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 250147.25
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 750441.75
$ws.Range("N23").Value = -750911.75

$ws.Range("H25").Value = 283.33334
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 250
$ws.Range("K25").Value = 900
$ws.Range("L25").Value = 750
$ws.Range("M25").Value = -731
$ws.Range("N25").Value = -1088

$ws.Range("H30").Value = 283.33334
$ws.Range("I30").Value = 300
$ws.Range("J30").Value = 250
$ws.Range("K30").Value = 900
$ws.Range("L30").Value = 750
$ws.Range("M30").Value = -798
$ws.Range("N30").Value = -954

$ws.Range("H61").Value = 27.25
$ws.Range("I61").Value = 24.666666
$ws.Range("J61").Value = 35
$ws.Range("K61").Value = 73.99999800000001
$ws.Range("L61").Value = 105
$ws.Range("M61").Value = 141.000002
$ws.Range("N61").Value = -535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 320.9
$ws.Range("I2").Value = 79.166664
$ws.Range("J2").Value = 683.5
$ws.Range("K2").Value = 79.166664
$ws.Range("L2").Value = 683.5
$ws.Range("M2").Value = 33.833336
$ws.Range("N2").Value = -909.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2450
$ws.Range("I22").Value = 1850
$ws.Range("J22").Value = 2750
$ws.Range("K22").Value = 1850
$ws.Range("L22").Value = 2750
$ws.Range("M22").Value = -1555
$ws.Range("N22").Value = -3340

$ws.Range("H27").Value = 2450
$ws.Range("I27").Value = 1850
$ws.Range("J27").Value = 2750
$ws.Range("K27").Value = 1850
$ws.Range("L27").Value = 2750
$ws.Range("M27").Value = -1743
$ws.Range("N27").Value = -2964

$ws.Range("H40").Value = 4500.8
$ws.Range("I40").Value = 3501.3333
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 3501.3333
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -3365.3333
$ws.Range("N40").Value = -6272

$ws.Range("H46").Value = 2844.6924
$ws.Range("I46").Value = 3001
$ws.Range("J46").Value = 2831.6667
$ws.Range("K46").Value = 3001
$ws.Range("L46").Value = 2831.6667
$ws.Range("M46").Value = -2813
$ws.Range("N46").Value = -3207.6667

$ws.Range("H68").Value = 3823.6667
$ws.Range("I68").Value = 3234.75
$ws.Range("J68").Value = 5001.5
$ws.Range("K68").Value = 3234.75
$ws.Range("L68").Value = 5001.5
$ws.Range("M68").Value = -2485.75
$ws.Range("N68").Value = -6499.5

$ws.Range("H71").Value = 3823.6667
$ws.Range("I71").Value = 3234.75
$ws.Range("J71").Value = 5001.5
$ws.Range("K71").Value = 16173.75
$ws.Range("L71").Value = 25007.5
$ws.Range("M71").Value = -12429.75
$ws.Range("N71").Value = -32495.5

$ws.Range("H132").Value = 4250
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -7970
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4010.182
$ws.Range("I62").Value = 2803.7778
$ws.Range("J62").Value = 4845.385
$ws.Range("K62").Value = 2803.7778
$ws.Range("L62").Value = 4845.385
$ws.Range("M62").Value = -2179.7778
$ws.Range("N62").Value = -6093.385

$ws.Range("H65").Value = 4010.182
$ws.Range("I65").Value = 2803.7778
$ws.Range("J65").Value = 4845.385
$ws.Range("K65").Value = 14018.889
$ws.Range("L65").Value = 24226.925
$ws.Range("M65").Value = -10898.889
$ws.Range("N65").Value = -30466.925

$ws.Range("H136").Value = 2356.2942
$ws.Range("I136").Value = 1326.2727
$ws.Range("J136").Value = 4244.6665
$ws.Range("K136").Value = 3978.8181
$ws.Range("L136").Value = 12733.9995
$ws.Range("M136").Value = -1428.8181
$ws.Range("N136").Value = -17833.9995
